{"js": "// Rewrite the (single) opening paragraph:\n//   \"this\" + \" is a test document.  Created in dev0.\"\n// becomes two paragraphs:\n//   \"this is a \" + \"test document.  Created in master.\"\n//   \"Edits from dev0.\"  (carrying the _GoBack bookmark)\n// and the stray grammar-check markers (<w:proofErr/>) around the old\n// \"this\" / \" is a ...\" run split are dropped, since they described a\n// sentence boundary that no longer exists after the rewrite.\nconst body = context.document.body;\n\nconst newBodyOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">this is a </w:t></w:r>\n            <w:r><w:t>test document.  Created in master.</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:t>Edits from dev0.</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nbody.insertOoxml(newBodyOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Rewrite the (single) opening paragraph:\n#   \"this\" + \" is a test document.  Created in dev0.\"\n# becomes two paragraphs:\n#   \"this is a \" + \"test document.  Created in master.\"\n#   \"Edits from dev0.\"  (carrying the _GoBack bookmark)\n# and the stray grammar-check markers (<w:proofErr/>) around the old\n# \"this\" / \" is a ...\" run split are dropped, since they described a\n# sentence boundary that no longer exists after the rewrite.\n$d = $word.ActiveDocument\n\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">this is a </w:t></w:r>\n            <w:r><w:t>test document.  Created in master.</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:t>Edits from dev0.</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$d.Content.InsertXML($xml)\n"}
